$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 2.9
$ws.Range("G3").Value = 4
$ws.Range("I3").Value = 2.74
$ws.Range("P3").Value = 1.73
$ws.Range("Q3").Value = 1.92
$ws.Range("V3").Value = 1.57
$ws.Range("F4").Value = 1.62
$ws.Range("G4").Value = 1.9
$ws.Range("H4").Value = 2.1
$ws.Range("I4").Value = 11
$ws.Range("P4").Value = 1.52
$ws.Range("P6").Value = 2.3
$ws.Range("G7").Value = 100
$ws.Range("H7").Value = 1.33
$ws.Range("Q7").Value = 1.87
$ws.Range("J8").Value = 3.2
$ws.Range("F9").Value = 2.84
$ws.Range("J9").Value = 2.64
$ws.Range("P9").Value = 2.04
$ws.Range("F10").Value = 1.69
$ws.Range("G10").Value = 1.99
$ws.Range("J10").Value = 3.35
$ws.Range("K10").Value = 6.2
$ws.Range("N10").Value = 1.72
$ws.Range("Q10").Value = 1.96
$ws.Range("S10").Value = 3.35
$ws.Range("V10").Value = 1.16
$ws.Range("W10").Value = 2
$ws.Range("F12").Value = 2.66
$ws.Range("H12").Value = 2.16
$ws.Range("I12").Value = 2.74
$ws.Range("J12").Value = 2.78
$ws.Range("K12").Value = 6.6
$ws.Range("P12").Value = 2.26
$ws.Range("Q12").Value = 1.53
$ws.Range("F13").Value = 3
$ws.Range("H13").Value = 1.89
$ws.Range("K13").Value = 7.6
$ws.Range("P13").Value = 3.15
$ws.Range("M17").Value = 1.12
$ws.Range("N17").Value = 2.26
$ws.Range("O17").Value = 1.52
$ws.Range("P17").Value = 1.43
$ws.Range("Q17").Value = 2.84
$ws.Range("G18").Value = 2.06
$ws.Range("N18").Value = 3.35
$ws.Range("F19").Value = 6.4
$ws.Range("G19").Value = 10.5
$ws.Range("I19").Value = 1.71
$ws.Range("J19").Value = 3.85
$ws.Range("K19").Value = 5.8
$ws.Range("P19").Value = 1.89
$ws.Range("Q19").Value = 1.74
$ws.Range("G20").Value = 2.28
$ws.Range("H20").Value = 3.5
$ws.Range("I20").Value = 3.9
$ws.Range("K20").Value = 3.9
$ws.Range("P20").Value = 1.96
$ws.Range("T20").Value = 1.72
$ws.Range("V20").Value = 1.34
$ws.Range("W20").Value = 1.79
$ws.Range("AD20").Value = 18.5
$ws.Range("AI20").Value = 60
$ws.Range("AO20").Value = 50
